$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.790.65'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '3.339.95'
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'580.19"
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = "'175.88"
$ws.Range('E6').Value = '  -1.74%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('D9').Value = '3.336.55'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('E10').Value = '  +2.27%  '
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').Value = "'46.30"
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').Value = "'708.95"
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').Value = '3.875.86'
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').Value = "'8.44"
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('D17').Value = '67.806.32'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('D19').Value = '3.340.33'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').Value = "'17.41"
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = "'11.01"
$ws.Range('E21').Value = '  +2.00%  '
$ws.Range('D22').Value = "'0.895"
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('E23').Value = '  +3.75%  '
$ws.Range('D24').Value = "'16.98"
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('D25').Value = "'98.67"
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').Value = "'3.88"
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('D28').Value = "'9.48"
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = "'33.32"
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').Value = "'8.54"
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('D31').Value = "'7.11"
$ws.Range('E31').Value = '  +5.34%  '
$ws.Range('D32').Value = "'571.79"
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').Value = "'10.99"
$ws.Range('E34').Value = '  +1.03%  '
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = "'57.37"
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '3.705.23'
$ws.Range('E37').Value = '  -3.98%  '
$ws.Range('D38').Value = "'3.36"
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('D39').Value = "'34.17"
$ws.Range('E39').Value = '  +6.18%  '
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = "'2.66"
$ws.Range('E41').Value = '  +1.19%  '
$ws.Range('D42').Value = "'3.19"
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').Value = '0.0₃0678'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('E45').Value = '  -3.18%  '
$ws.Range('D46').Value = "'0.0406"
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('E47').Value = '  +6.28%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('E50').Value = '  -5.69%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').Value = "'2.71"
$ws.Range('E51').Value = '  +16.65%  '
